$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Header / summary updates ---
$ws.Range("E11").Value = 246820
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 1

# --- Insert two additional data rows (table goes from 2 rows to 4 rows) ---
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

# Copy formatting from the still-intact first data row (row 16) onto the two
# newly-inserted blank rows (17 and 18) so borders/fonts/fills match.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 16: MARLENE DEL CARMEN VITAL ACOSTA ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "42206362"
$ws.Range("D16").Value = "MARLENE DEL CARMEN VITAL ACOSTA"
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 76000
$ws.Range("G16").Value = 1900000

# --- Row 17: HASMED MOISES CASTRO VITAL ---
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143356633"
$ws.Range("D17").Value = "HASMED MOISES CASTRO VITAL"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# --- Row 18: DANIEL ANDRES CASTRO VITAL ---
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1007938749"
$ws.Range("D18").Value = "DANIEL ANDRES CASTRO VITAL"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# --- Row 19 (previously row 17): MARNELLY ACOSTA POLANCO, period updated ---
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1002280059"
$ws.Range("D19").Value = "MARNELLY ACOSTA POLANCO"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500
